# estudoNN_TP.xlsx — "Completei a primeira tabela e mudei nn_train para
# fazer a media dos 10 automaticamente"
#
# Fills in the remaining rows (Conf3/Conf4) of the first results table
# ("A função de treino influencia o desempenho?") with the missing
# hidden-layer / neuron-count / activation-function configuration and the
# measured accuracy results, matching the style already used by the
# Conf1/Conf2 rows above them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IRIS")

# --- Conf1 (row 7) / Conf2 (row 8): accuracies now stored as real numbers
$ws.Range("H7").Value = 98.89
$ws.Range("I7").Value = 97.88

$ws.Range("H8").Value = 99.26
$ws.Range("I8").Value = 97.89

# --- Conf3 (row 9): complete the row
$ws.Range("B9").Value = 3
$ws.Range("C9").Value = "5,10,5"
$ws.Range("D9").Value = "tansig, tansig,tansig, purelin"
$ws.Range("H9").Value = 98.78
$ws.Range("I9").Value = 96.5

# Match the bold red "highlighted" look used for B7:C8
$ws.Range("B9:C9").Font.Bold = $true
$ws.Range("B9:C9").Font.Color = 255

# --- Conf4 (row 10): complete the row
$ws.Range("B10").Value = 3
$ws.Range("C10").Value = "10,10,10"
$ws.Range("D10").Value = "tansig, tansig,tansig, purelin"
$ws.Range("H10").Value = 98.834000000000003
$ws.Range("I10").Value = 96.61

$ws.Range("B10:C10").Font.Bold = $true
$ws.Range("B10:C10").Font.Color = 255

# Move the active selection, matching the author's final cursor position
$null = $ws.Range("K10").Select()
